# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.885.98'
$ws.Range('E2').Value = '  +0.16%  '

# Row 3
$ws.Range('D3').Value = '2.626.94'
$ws.Range('E3').Value = '  -1.21%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.09'
$ws.Range('E5').Value = '  -1.29%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.16'
$ws.Range('E6').Value = '  +1.41%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('E8').Value = '  -0.23%  '

# Row 9
$ws.Range('E9').Value = '  +0.62%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.65'
$ws.Range('E10').Value = '  +0.82%  '

# Row 11
$ws.Range('E11').Value = '  +3.42%  '

# Row 12
$ws.Range('E12').Value = '  -1.09%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.67'
$ws.Range('E13').Value = '  +0.03%  '

# Row 14
$ws.Range('D14').Value = '3.098.04'
$ws.Range('E14').Value = '  -1.19%  '

# Row 15
$ws.Range('D15').Value = '63.718.90'
$ws.Range('E15').Value = '  +0.10%  '

# Row 16
$ws.Range('E16').Value = '  +2.17%  '

# Row 17
$ws.Range('D17').Value = '2.616.61'
$ws.Range('E17').Value = '  -1.09%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.32'
$ws.Range('E18').Value = '  +7.02%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.64'
$ws.Range('E19').Value = '  +1.95%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '349.99'
$ws.Range('E20').Value = '  +2.00%  '

# Row 21
$ws.Range('E21').Value = '  -1.29%  '

# Row 22
$ws.Range('E22').Value = '  -0.26%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.70'
$ws.Range('E23').Value = '  +2.15%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.36'
$ws.Range('E24').Value = '  -0.91%  '

# Row 25
$ws.Range('E25').Value = '  +12.37%  '

# Row 26
$ws.Range('E26').Value = '  +0.46%  '

# Row 27
$ws.Range('E27').Value = '  -2.42%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '561.63'
$ws.Range('E28').Value = '  +1.01%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.21'
$ws.Range('E29').Value = '  +4.45%  '

# Row 30
$ws.Range('E30').Value = '  +0.16%  '

# Row 31
$ws.Range('E31').Value = '  +0.12%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.05'
$ws.Range('E32').Value = '  -0.21%  '

# Row 33
$ws.Range('D33').Value = '0.0₃0845'
$ws.Range('E33').Value = '  +2.51%  '

# Row 34
$ws.Range('E34').Value = '  -1.49%  '

# Row 35
$ws.Range('E35').Value = '  +0.34%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '169.41'
$ws.Range('E36').Value = '  +1.10%  '

# Row 37
$ws.Range('E37').Value = '  +0.52%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.03%  '

# Row 39
$ws.Range('E39').Value = '  +0.65%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.38'
$ws.Range('E40').Value = '  +1.06%  '

# Row 41
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.06%  '

# Row 42
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '170.13'
$ws.Range('E42').Value = '  +0.57%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.86'
$ws.Range('E43').Value = '  -0.27%  '

# Row 44
$ws.Range('E44').Value = '  +3.51%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0599'
$ws.Range('E45').Value = '  +3.65%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.48'
$ws.Range('E46').Value = '  -4.89%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.630'
$ws.Range('E47').Value = '  -0.15%  '

# Row 48
$ws.Range('E48').Value = '  +0.02%  '

# Row 49
$ws.Range('E49').Value = '  +5.17%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0969'
$ws.Range('E50').Value = '  +0.36%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.25'
$ws.Range('E51').Value = '  +1.90%  '
